$d = $word.ActiveDocument

# 1) Insert a new bullet paragraph "Added the "Watch Ads" menu (not yet
#    functional)." right before the "Added Max Gold limit..." bullet,
#    matching its list formatting.
$r = $d.Content
$r.Find.Execute("Added Max Gold limit to Marketplace transactions.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $r.Paragraphs(1)
$target.Range.InsertParagraphBefore()
$target.Range.Text = "Added the " + [char]0x201C + "Watch Ads" + [char]0x201D + " menu (not yet functional)."

# 2) Reword the "Slightly increased text size in Quest briefings." bullet.
$d.Content.Find.Execute("Slightly increased text size in Quest briefings.", $true, $false, $false, $false, $false, $true, 1, $false, "Some visual adjustments to many menus.", 2)
